# Generate Report for Handback
# Adds a new handback record (row 8, "Latest Handback" columns I/J/K + Error Detail
# column P) for the e963eb02-9bdf-46fa-b9a7-262f9c8bbbd7 file on both the zh-cn and
# de-de localization-status worksheets, together with a hyperlink on the new
# "Latest Handback File" cell and a widened "Error Detail" column.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b16452da664330d2efb66f35deafa91c8efc49e/e2e/e963eb02-9bdf-46fa-b9a7-262f9c8bbbd7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2c3c827dc13070aa3e3a6123c4cbbe233d3b709/e2e/e963eb02-9bdf-46fa-b9a7-262f9c8bbbd7.md."

function Update-HandbackRow {
    param(
        [string]$SheetName,
        [string]$TargetFile,
        [string]$HandbackDateTime,
        [string]$HyperlinkTarget
    )

    $ws = $wb.Worksheets.Item($SheetName)
    Write-Host "Updating sheet:" $SheetName "target:" $TargetFile "dt:" $HandbackDateTime "link:" $HyperlinkTarget

    # Latest Handback File (display name, hyperlinked to the handback .md)
    $ws.Range("I8").Value = "e963eb02-9bdf-46fa-b9a7-262f9c8bbbd7.md"
    $ws.Hyperlinks.Add($ws.Range("I8"), $HyperlinkTarget, "", "", "e963eb02-9bdf-46fa-b9a7-262f9c8bbbd7.md") | Out-Null
    $ws.Range("I8").Font.Underline = 2
    $ws.Range("I8").Font.Color = 15570276

    # Latest Target File (the generated xlf for this handback)
    $ws.Range("J8").Value = $TargetFile

    # Latest Handback DateTime
    $ws.Range("K8").Value = $HandbackDateTime

    # Error Detail
    $ws.Range("P8").Value = $errorDetail

    # Error Detail column needs to be wide enough to read the message
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

Update-HandbackRow `
    "zh-cn" `
    "e963eb02-9bdf-46fa-b9a7-262f9c8bbbd7.428e36d09ab1195149301d3aba4676147da8572d.zh-cn.xlf" `
    "2016-08-20 08:51:25" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/23458073678df38eb3f7fa6b421d5480117f608b/e2e/e963eb02-9bdf-46fa-b9a7-262f9c8bbbd7.md"

Update-HandbackRow `
    "de-de" `
    "e963eb02-9bdf-46fa-b9a7-262f9c8bbbd7.428e36d09ab1195149301d3aba4676147da8572d.de-de.xlf" `
    "2016-08-20 08:51:32" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/81bd31b2a832704f2ef8133ffad3685011d30cb5/e2e/e963eb02-9bdf-46fa-b9a7-262f9c8bbbd7.md"
